$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ivy"
$ws.Range("D2").Value = "Shake with ice for a long time. Strain into a chilled g{cocktail glass}."
$ws.Range("B2").Value = "dry gin, green chartreuse, dry vermouth, absinthe, orange bitters"
$ws.Range("C2").Value = "q{1.5} u{oz} dry gin | q{1/2} u{oz} dry vermouth | q{1/2} u{oz} green Chartreuse | q{1} u{tsp} absinthe | q{1-2} u{dashes} orange bitters"

$ws.Range("B1").Select()
